$d = $word.ActiveDocument

# 1) Update the summary/title table (Table 1) entry from the "TFS 25016" text to the
#    new "TFS 25205" text.
$t1 = $d.Tables.Item(1)
$t1.Cell(1, 2).Range.Text = "25205 – New Submission: display log name with success message for single log submission."

# 2) Append a new row to the change log table (Table 2) documenting the 25205 change.
$t2 = $d.Tables.Item(2)
$t2.Rows.Add() | Out-Null
$newRowIndex = $t2.Rows.Count
$t2.Cell($newRowIndex, 1).Range.Text = "08/15/2022"
$t2.Cell($newRowIndex, 2).Range.Text = "TFS 25205 – New Submission: display log name with success message for single log submission."
$t2.Cell($newRowIndex, 3).Range.Text = "Lili Huang"

# 3) Bump the changeset number (51750 -> 51847) referenced in the implementation steps
#    table (Table 3). We edit only the individual digit characters (rather than
#    replacing the whole phrase) so the surrounding formatting - the plain
#    "Changeset " label, the large bold changeset number, and the small plain
#    trailing ";" - stays exactly as it was.
$searchRange = $d.Content
$null = $searchRange.Find.Execute("Changeset 51750;", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
$phraseStart = $searchRange.Start
# "Changeset " is 10 characters, so the digits "51750" start right after it.
$digitsStart = $phraseStart + 10
# Remove the trailing "0" (5th digit) first, then the "5" (4th digit), so offsets
# of the earlier characters stay valid; finally turn the 3rd digit ("7") into "847".
$d.Range($digitsStart + 4, $digitsStart + 5).Text = ""
$d.Range($digitsStart + 3, $digitsStart + 4).Text = ""
$d.Range($digitsStart + 2, $digitsStart + 3).Text = "847"
